# Daily attendance processing - 2025-12-21 23:27:33
# Reorders the "Recorded By" email list (column G) for several rows on the
# "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value  = "servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System, gehanadel@med.asu.edu.eg"
$ws.Range("G3").Value  = "majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G4").Value  = "majorelle.magdy@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
$ws.Range("G6").Value  = "majorelle.magdy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, manar.montaser@med.asu.edu.eg"
$ws.Range("G7").Value  = "Amera.a.saad@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg"
$ws.Range("G8").Value  = "AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg"
$ws.Range("G11").Value = "Safa.hany@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
$ws.Range("G12").Value = "Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg"
$ws.Range("G13").Value = "esraa.mostafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg"
$ws.Range("G15").Value = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("G19").Value = "mariam.youssif.std@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("G24").Value = "Sarah.Mahdy@med.asu.edu.eg, youstina.gamil@med.asu.edu.eg"
$ws.Range("G28").Value = "maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"
$ws.Range("G30").Value = "wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
